# Apply the "update database and change read_price algorithm" change:
# - shift the fiscal-period columns (D..H) one period to the left and
#   append the newest period's figures in column H
# - update the period-label row (8) and the publish-date row (9) accordingly
# - fix up the one stray "-" placeholder (row 15, col D) into a real 0

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: "12 ماهه منتهی به ..." period headers (D8:H8) ---
$ws.Cells.Item(8, 4).Value = "12 ماهه منتهی به 1397/12"
$ws.Cells.Item(8, 5).Value = "12 ماهه منتهی به 1398/12"
$ws.Cells.Item(8, 6).Value = "12 ماهه منتهی به 1399/12"
$ws.Cells.Item(8, 7).Value = "12 ماهه منتهی به 1400/12"
$ws.Cells.Item(8, 8).Value = "12 ماهه منتهی به 1401/12"

# --- Row 9: "تاریخ انتشار" publish dates (D9:H9) ---
$ws.Cells.Item(9, 4).Value = "1399-02-21 (8)"
$ws.Cells.Item(9, 5).Value = "1400-02-29 (9)"
$ws.Cells.Item(9, 6).Value = "1401-03-04 (8)"
$ws.Cells.Item(9, 7).Value = "1402-02-30 (8)"
$ws.Cells.Item(9, 8).Value = "1402-02-30 (2)"

# --- Financial data rows: shift values left by one period and add the new period ---
$dataRows = @{
    11 = @(4371785, 7801585, 10395584, 21171738, 25320705)
    12 = @(-2452229, -3794403, -5016199, -9937618, -13409448)
    13 = @(1919556, 4007182, 5379385, 11234120, 11911257)
    14 = @(-217408, -317404, -403929, -839072, -971858)
    15 = @(0, 0, 0, 0, 0)
    16 = @(42505, 9985, -17261, -52392, -9710)
    17 = @(1744653, 3699763, 4958195, 10342656, 10929689)
    18 = @(-548697, -373314, -531568, -959126, -1783257)
    19 = @(77529, 51002, 191129, 62162, 265292)
    20 = @(1273485, 3377451, 4617756, 9445692, 9411724)
    21 = @(-225543, -588307, -642189, -1360326, -1498403)
    22 = @(1047942, 2789144, 3975567, 8085366, 7913321)
    24 = @(1047942, 2789144, 3975567, 8085366, 7913321)
    25 = @(2329, 6198, 3534, 7187, 3517)
    26 = @(450000, 450000, 1125000, 1125000, 2250000)
    27 = @(466, 1240, 1767, 3594, 3517)
}

foreach ($r in $dataRows.Keys) {
    $values = $dataRows[$r]
    for ($i = 0; $i -lt 5; $i++) {
        $ws.Cells.Item($r, 4 + $i).Value = $values[$i]
    }
}
